# Language table update: thawing explanation and details
#
# - sanitaryEnd_1 (row 104, column B) text is tightened: dropped the
#   redundant "the food we need" clause.
# - 27 new localization rows are appended (rows 105-131) introducing the
#   TCS-food explanation, the produce-washing intro, the three thawing
#   methods (fridge / running water / microwave) with their VoiceDuration
#   numbers in column C, the thaw exercise dialogue, and the thermometer
#   calibration intro.
#
# Cell writes below are ordered to match the original authoring order
# (mostly row-by-row, with a couple of column-batched spots) so the
# resulting shared-string table lines up with the target workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Tighten the existing sanitaryEnd_1 value.
$ws.Cells.Item(104, 2).Value = "Alright, we are good to go! Let's go ahead and start prepping for the day."
$ws.Cells.Item(105, 1).Value = "tcsIntro_1"
$ws.Cells.Item(105, 2).Value = "First thing we need to do is determine which foods are potentially hazardous. The ones that allow pathogen growth."
$ws.Cells.Item(106, 1).Value = "tcsIntro_2"
$ws.Cells.Item(106, 2).Value = "The best way to deal with these foods is through proper time management, and temperature control."
$ws.Cells.Item(107, 1).Value = "tcsIntro_3"
$ws.Cells.Item(107, 2).Value = "This is known as Time and Temperature Control for Safety, or in short: TCS food."
$ws.Cells.Item(108, 1).Value = "tcsIntro_4"
$ws.Cells.Item(108, 2).Value = "Most of these foods are easily identified through common sense. So, let’s see if you can identify which ones are TCS or non-TCS food."
$ws.Cells.Item(109, 1).Value = "tcsPost_1"
$ws.Cells.Item(109, 2).Value = "Excellent! Now that that’s sorted out, let me show you a comprehensive list for TCS food."
$ws.Cells.Item(110, 1).Value = "produceIntro_1"
$ws.Cells.Item(110, 2).Value = "Since the majority of foods already carry germs, steps must be taken to minimize the potential danger of these foods possess."
$ws.Cells.Item(111, 1).Value = "produceIntro_2"
$ws.Cells.Item(111, 2).Value = "As such, all produce needs to be thoroughly washed to remove germs, and wash away insecticides."
$ws.Cells.Item(112, 1).Value = "produceIntro_3"
$ws.Cells.Item(112, 2).Value = "Here are some useful tips for washing produce."
$ws.Cells.Item(113, 1).Value = "thawMethods_title"
$ws.Cells.Item(113, 2).Value = "Thawing Methods"
$ws.Cells.Item(114, 1).Value = "thawMethods_fridge_title"
$ws.Cells.Item(114, 2).Value = "Refridgerator"
$ws.Cells.Item(114, 3).Value = 1.5
$ws.Cells.Item(115, 1).Value = "thawMethods_fridge_desc"
$ws.Cells.Item(115, 2).Value = "Thaw food in a refrigerator at 41° F (5° C) or lower to keep dangerous microorganisms from growing. Plan ahead when thawing large items such as turkeys - they can take several days to thaw."
$ws.Cells.Item(115, 3).Value = 8
$ws.Cells.Item(116, 1).Value = "thawMethods_water_title"
$ws.Cells.Item(117, 1).Value = "thawMethods_water_desc"
$ws.Cells.Item(116, 2).Value = "Running Water"
$ws.Cells.Item(116, 3).Value = 2
$ws.Cells.Item(117, 2).Value = "Thaw food submerged under running water at a temperature of 70° F (21° C) or lower. The water flow must be strong enough to wash food particles into the overflow drain."
$ws.Cells.Item(117, 3).Value = 8
$ws.Cells.Item(118, 1).Value = "thawMethods_microwave_title"
$ws.Cells.Item(119, 1).Value = "thawMethods_microwave_desc"
$ws.Cells.Item(118, 2).Value = "Microwave"
$ws.Cells.Item(118, 3).Value = 1
$ws.Cells.Item(119, 2).Value = "You can safely thaw food in a microwave oven if the food will be cooked immediately. Large items such as roasts or turkeys may not thaw well in a microwave."
$ws.Cells.Item(119, 3).Value = 8
$ws.Cells.Item(120, 1).Value = "thawIntro_1"
$ws.Cells.Item(120, 2).Value = "Now we will talk about proper thawing of food."
$ws.Cells.Item(121, 1).Value = "thawIntro_2"
$ws.Cells.Item(122, 1).Value = "thawIntro_3"
$ws.Cells.Item(121, 2).Value = "There are several ways to thaw food. Each method varies by how long it takes. So, plan accordingly!"
$ws.Cells.Item(122, 2).Value = "Here are the three methods commonly used."
$ws.Cells.Item(123, 1).Value = "thawFaucet_1"
$ws.Cells.Item(124, 1).Value = "thawFaucet_2"
$ws.Cells.Item(123, 2).Value = "For this exercise, we will be using the running water method."
$ws.Cells.Item(124, 2).Value = "Let’s start by turning the faucet on."
$ws.Cells.Item(125, 1).Value = "thawMeat_1"
$ws.Cells.Item(125, 2).Value = "Great, now grab the pack of meat and put it under the running water."
$ws.Cells.Item(126, 1).Value = "thawWait_1"
$ws.Cells.Item(126, 2).Value = "Excellent! Now just wait till the meat has properly thawed…It should take about half an hour to an hour."
$ws.Cells.Item(127, 2).Value = "Alright, well this will obviously take a while, so let’s compress time to speed things up!"
$ws.Cells.Item(127, 1).Value = "thawWait_2"
$ws.Cells.Item(128, 1).Value = "thawEnd_1"
$ws.Cells.Item(128, 2).Value = "Great! The meat has been properly thawed and is ready to be cooked."
$ws.Cells.Item(129, 1).Value = "thawEnd_2"
$ws.Cells.Item(129, 2).Value = "Remember that once food has been thawed, it must be cooked right away! Don’t put any of them back in the freezer!"
$ws.Cells.Item(130, 1).Value = "calibrateIntro_1"
$ws.Cells.Item(130, 2).Value = "Let's now talk about properly calibrating a thermometer."
$ws.Cells.Item(131, 1).Value = "calibrateIntro_2"
$ws.Cells.Item(131, 2).Value = "Before using a thermometer, make sure it is calibrated properly!"

# Keep the selection in sync with the new bottom of the table.
$ws.Range("B132").Select()
